# Update a handful of computed tax-support share values on "Sheet 1"
# after re-running the data prepare & render pipeline with final data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet 1")

# Row 2 (country/category row for A2 = key 14)
$ws.Range("B2").Value = 0.556007222541788
$ws.Range("K2").Value = 0.438282499020205
$ws.Range("L2").Value = 0.595756191953926
$ws.Range("N2").Value = 0.508231644030169

# Row 3 (country/category row for A3 = key 15)
$ws.Range("B3").Value = 0.50312291438834
$ws.Range("K3").Value = 0.351050882580874
$ws.Range("L3").Value = 0.604399737467109
$ws.Range("N3").Value = 0.447300643788012
